$d = $word.ActiveDocument
$sec = $d.Sections.First

# Pearson logo pictures live in the footers (docPr/cNvPr name "image2.png" -> "image1.png").
# Headers(1)/Footers(1) = primary ("default") header/footer; Headers(2)/Footers(2) = first-page
# header/footer. Selecting the shape's own Range before touching Selection.InlineShapes(1).Name
# keeps the rename anchored on the picture actually being edited.

$footerPrimary = $sec.Footers.Item(1)
$ishp = $footerPrimary.Range.InlineShapes.Item(1)
$ishp.Range.Select
$word.Selection.InlineShapes.Item(1).Name = "image1.png"

$footerFirst = $sec.Footers.Item(2)
$ishp = $footerFirst.Range.InlineShapes.Item(1)
$ishp.Range.Select
$word.Selection.InlineShapes.Item(1).Name = "image1.png"

# BTec logo pictures live in the headers ("image1.jpg" -> "image2.jpg").
$headerPrimary = $sec.Headers.Item(1)
$ishp = $headerPrimary.Range.InlineShapes.Item(1)
$ishp.Range.Select
$word.Selection.InlineShapes.Item(1).Name = "image2.jpg"

$headerFirst = $sec.Headers.Item(2)
$ishp = $headerFirst.Range.InlineShapes.Item(1)
$ishp.Range.Select
$word.Selection.InlineShapes.Item(1).Name = "image2.jpg"
